$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the format (font/border/alignment)
# used by the other header cells (e.g. G1) so H1 gets the same style record.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# Add the new data value in H2 (numeric 0, default/unstyled like the other data cells)
$ws.Range("H2").Value = 0
